$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.378.79'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '2.774.74'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '352.47'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.80'
$ws.Range('E6').Value = '  -1.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.548'
$ws.Range('E7').Value = '  -2.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.588'
$ws.Range('E9').Value = '  -1.20%  '
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('E11').Value = '  +3.06%  '
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('E13').Value = '  +3.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.60'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = '3.210.83'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').Value = '2.766.08'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.921'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').Value = '51.389.77'
$ws.Range('E18').Value = '  -0.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.61'
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('E20').Value = '  -1.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.10'
$ws.Range('E21').Value = '  +0.80%  '
$ws.Range('D22').Value = '0.0₃0961'
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.85'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '265.34'
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -2.13%  '
$ws.Range('E28').Value = '  +12.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.21'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.20'
$ws.Range('E30').Value = '  +1.97%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.28'
$ws.Range('E31').Value = '  +6.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.21'
$ws.Range('E32').Value = '  +9.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '51.91'
$ws.Range('E33').Value = '  +0.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0455'
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.56'
$ws.Range('E35').Value = '  +6.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0826'
$ws.Range('E36').Value = '  -2.37%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.52'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '120.35'
$ws.Range('E43').Value = '  -0.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.00'
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('D46').Value = '2.103.93'
$ws.Range('E46').Value = '  +1.85%  '
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('E48').Value = '  +5.08%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.41'
$ws.Range('E49').Value = '  -4.86%  '
$ws.Range('B50').Value = 'SEI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.904'
$ws.Range('E50').Value = '  -3.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.31'
$ws.Range('E51').Value = '  +7.83%  '
